# Apply the latest cryptos snapshot values scraped on Wed Mar 15 17:16:10 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.363.68'
$ws.Range("E2").Value = '  -6.09%  '
$ws.Range("D3").Value = '1.636.53'
$ws.Range("E3").Value = '  -7.51%  '
$ws.Range("D4").Value = '''1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.61%  '
$ws.Range("D5").Value = '''1.003'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.48%  '
$ws.Range("D6").Value = '''305.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.91%  '
$ws.Range("D7").Value = '''0.3616'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.69%  '
$ws.Range("D8").Value = '''46.92'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.49%  '
$ws.Range("D9").Value = '''0.3219'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -11.08%  '
$ws.Range("D10").Value = '''1.103'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -10.76%  '
$ws.Range("D11").Value = '''0.06892'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -10.45%  '
$ws.Range("D12").Value = '''1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").Value = '''5.902'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -9.23%  '
$ws.Range("D14").Value = '''19.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -12.51%  '
$ws.Range("D15").Value = '1.641.08'
$ws.Range("E15").Value = '  -7.21%  '
$ws.Range("D16").Value = '''6.501'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.53%  '
$ws.Range("D17").Value = '''0.00001038'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -10.47%  '
$ws.Range("E18").Value = '  -4.21%  '
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D20").Value = '''77.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -11.99%  '
$ws.Range("D21").Value = '''15.73'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -11.49%  '
$ws.Range("D22").Value = '''5.874'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -9.93%  '
$ws.Range("D23").Value = '''11.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.80%  '
$ws.Range("D24").Value = '24.364.34'
$ws.Range("E24").Value = '  -5.80%  '
$ws.Range("D25").Value = '''2.413'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("D26").Value = '''2.368'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -19.42%  '
$ws.Range("D27").Value = '''144.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.10%  '
$ws.Range("D28").Value = '''18.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -10.05%  '
$ws.Range("D29").Value = '1.810.23'
$ws.Range("E29").Value = '  -7.85%  '
$ws.Range("D30").Value = '''124.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.53%  '
$ws.Range("D31").Value = '''1.087'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.42%  '
$ws.Range("D32").Value = '''4.055'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.24%  '
$ws.Range("D33").Value = '''5.672'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -21.63%  '
$ws.Range("D34").Value = '''0.08356'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.12%  '
$ws.Range("D35").Value = '''1.668'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.33%  '
$ws.Range("D36").Value = '''12.33'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -13.96%  '
$ws.Range("D37").Value = '''5.100'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -11.04%  '
$ws.Range("D38").Value = '''0.05985'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -11.71%  '
$ws.Range("D39").Value = '''0.02209'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -11.71%  '
$ws.Range("D40").Value = '''1.201'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.45%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '''0.2035'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -10.18%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''8.114'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -13.90%  '
$ws.Range("D43").Value = '''1.003'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.53%  '
$ws.Range("D44").Value = '''0.5859'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -11.57%  '
$ws.Range("D45").Value = '''3.735'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.70%  '
$ws.Range("D46").Value = '''12.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -13.51%  '
$ws.Range("D47").Value = '''0.5535'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -13.24%  '
$ws.Range("D48").Value = '''122.13'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.24%  '
$ws.Range("D49").Value = '''1.914'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -12.25%  '
$ws.Range("D50").Value = '''0.06907'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.20%  '
$ws.Range("D51").Value = '''73.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.80%  '
